# Auto-generated edit script applying scheduled-runner price/profit updates
# to the Anima_Profits workbook (per-sheet Leve profit tables).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 8436.24
$ws.Range("I19").Value = 292.1
$ws.Range("J19").Value = 13865.667
$ws.Range("K19").Value = 292.1
$ws.Range("L19").Value = 13865.667
$ws.Range("M19").Value = -117.1
$ws.Range("N19").Value = -14215.667
$ws.Range("H51").Value = 3775
$ws.Range("I51").Value = 1500
$ws.Range("J51").Value = 4533.3335
$ws.Range("K51").Value = 1500
$ws.Range("L51").Value = 4533.3335
$ws.Range("M51").Value = -1016
$ws.Range("N51").Value = -5501.3335
$ws.Range("H74").Value = 2837.5
$ws.Range("I74").Value = 2000
$ws.Range("K74").Value = 2000
$ws.Range("M74").Value = -1064
$ws.Range("H77").Value = 2837.5
$ws.Range("I77").Value = 2000
$ws.Range("K77").Value = 10000
$ws.Range("M77").Value = -5320
$ws.Range("H121").Value = 1448.0435
$ws.Range("J121").Value = 1507.0454
$ws.Range("L121").Value = 4521.1362
$ws.Range("N121").Value = -8015.1362
$ws.Range("H137").Value = 3335624.5
$ws.Range("I137").Value = 5954116.5
$ws.Range("J137").Value = 2998.7273
$ws.Range("K137").Value = 17862349.5
$ws.Range("L137").Value = 8996.1819
$ws.Range("M137").Value = -17859799.5
$ws.Range("N137").Value = -14096.1819
$ws.Range("H138").Value = 4753.3447
$ws.Range("I138").Value = 6002.6
$ws.Range("J138").Value = 4493.0835
$ws.Range("K138").Value = 18007.8
$ws.Range("L138").Value = 13479.2505
$ws.Range("M138").Value = -12867.8
$ws.Range("N138").Value = -23759.2505

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H63").Value = 32829.22
$ws.Range("J63").Value = 4915
$ws.Range("L63").Value = 4915
$ws.Range("N63").Value = -6287
$ws.Range("H66").Value = 32829.22
$ws.Range("J66").Value = 4915
$ws.Range("L66").Value = 24575
$ws.Range("N66").Value = -31439
$ws.Range("H74").Value = 10205675
$ws.Range("I74").Value = 951
$ws.Range("J74").Value = 33336382
$ws.Range("K74").Value = 951
$ws.Range("L74").Value = 33336382
$ws.Range("M74").Value = -77
$ws.Range("N74").Value = -33338130
$ws.Range("H77").Value = 10205675
$ws.Range("I77").Value = 951
$ws.Range("J77").Value = 33336382
$ws.Range("K77").Value = 4755
$ws.Range("L77").Value = 166681910
$ws.Range("M77").Value = -387
$ws.Range("N77").Value = -166690646
$ws.Range("H132").Value = 1427353.6
$ws.Range("I132").Value = 2600.641
$ws.Range("J132").Value = 5131711.5
$ws.Range("K132").Value = 7801.923000000001
$ws.Range("L132").Value = 15395134.5
$ws.Range("M132").Value = -5271.923000000001
$ws.Range("N132").Value = -15400194.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 79232.73
$ws.Range("J132").Value = 79156
$ws.Range("L132").Value = 79156
$ws.Range("N132").Value = -89276

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6442.15
$ws.Range("I31").Value = 2233.389
$ws.Range("J31").Value = 9885.682000000001
$ws.Range("K31").Value = 2233.389
$ws.Range("L31").Value = 9885.682000000001
$ws.Range("M31").Value = -1938.389
$ws.Range("N31").Value = -10475.682
$ws.Range("H34").Value = 6442.15
$ws.Range("I34").Value = 2233.389
$ws.Range("J34").Value = 9885.682000000001
$ws.Range("K34").Value = 2233.389
$ws.Range("L34").Value = 9885.682000000001
$ws.Range("M34").Value = -2031.389
$ws.Range("N34").Value = -10289.682
$ws.Range("H132").Value = 14816766
$ws.Range("I132").Value = 20001720
$ws.Range("J132").Value = 8335574
$ws.Range("K132").Value = 60005160
$ws.Range("L132").Value = 25006722
$ws.Range("M132").Value = -60002630
$ws.Range("N132").Value = -25011782
$ws.Range("H134").Value = 11370983
$ws.Range("I134").Value = 13165033
$ws.Range("J134").Value = 8668.333000000001
$ws.Range("K134").Value = 39495099
$ws.Range("L134").Value = 26004.999
$ws.Range("M134").Value = -39492564
$ws.Range("N134").Value = -31074.999
$ws.Range("H140").Value = 64290
$ws.Range("J140").Value = 64290
$ws.Range("L140").Value = 64290
$ws.Range("N140").Value = -74650
$ws.Range("H141").Value = 74351.84
$ws.Range("J141").Value = 71429.78999999999
$ws.Range("L141").Value = 71429.78999999999
$ws.Range("N141").Value = -81789.78999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 695.37036
$ws.Range("I68").Value = 619.3333
$ws.Range("J68").Value = 961.5
$ws.Range("K68").Value = 1857.9999
$ws.Range("L68").Value = 2884.5
$ws.Range("M68").Value = -1046.9999
$ws.Range("N68").Value = -4506.5
$ws.Range("H71").Value = 695.37036
$ws.Range("I71").Value = 619.3333
$ws.Range("J71").Value = 961.5
$ws.Range("K71").Value = 5573.9997
$ws.Range("L71").Value = 8653.5
$ws.Range("M71").Value = -1517.9997
$ws.Range("N71").Value = -16765.5
$ws.Range("H113").Value = 702.05554
$ws.Range("I113").Value = 688.6842
$ws.Range("J113").Value = 717
$ws.Range("K113").Value = 2066.0526
$ws.Range("L113").Value = 2151
$ws.Range("M113").Value = 103.9474
$ws.Range("N113").Value = -6491
$ws.Range("H122").Value = 2938.698
$ws.Range("I122").Value = 552.25
$ws.Range("K122").Value = 4970.25
$ws.Range("M122").Value = -2520.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 28785.8
$ws.Range("I40").Value = 10000
$ws.Range("J40").Value = 33482.25
$ws.Range("K40").Value = 10000
$ws.Range("L40").Value = 33482.25
$ws.Range("M40").Value = -9851
$ws.Range("N40").Value = -33780.25
$ws.Range("H124").Value = 36473
$ws.Range("J124").Value = 36473
$ws.Range("L124").Value = 36473
$ws.Range("N124").Value = -46293
$ws.Range("H132").Value = 5835709.5
$ws.Range("I132").Value = 2218.6843
$ws.Range("J132").Value = 24308430
$ws.Range("K132").Value = 6656.0529
$ws.Range("L132").Value = 72925290
$ws.Range("M132").Value = -4126.0529
$ws.Range("N132").Value = -72930350
